$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 4133.3335
$ws.Range("I118").Value = 400
$ws.Range("J118").Value = 6000
$ws.Range("K118").Value = 1200
$ws.Range("L118").Value = 18000
$ws.Range("M118").Value = 457
$ws.Range("N118").Value = -21314
$ws.Range("H137").Value = 7408266.5
$ws.Range("I137").Value = 781.8571
$ws.Range("K137").Value = 2345.5713
$ws.Range("M137").Value = 204.4287000000004
$ws.Range("H139").Value = 139716
$ws.Range("J139").Value = 139716
$ws.Range("L139").Value = 139716
$ws.Range("N139").Value = -149996
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9260541
$ws.Range("I61").Value = 10418001
$ws.Range("J61").Value = 860.5
$ws.Range("K61").Value = 10418001
$ws.Range("L61").Value = 860.5
$ws.Range("M61").Value = -10417789
$ws.Range("N61").Value = -1284.5
$ws.Range("H74").Value = 10871299
$ws.Range("I74").Value = 15153289
$ws.Range("J74").Value = 1631.3846
$ws.Range("K74").Value = 15153289
$ws.Range("L74").Value = 1631.3846
$ws.Range("M74").Value = -15152415
$ws.Range("N74").Value = -3379.3846
$ws.Range("H77").Value = 10871299
$ws.Range("I77").Value = 15153289
$ws.Range("J77").Value = 1631.3846
$ws.Range("K77").Value = 75766445
$ws.Range("L77").Value = 8156.923000000001
$ws.Range("M77").Value = -75762077
$ws.Range("N77").Value = -16892.923
$ws.Range("H136").Value = 9260541
$ws.Range("I136").Value = 10418001
$ws.Range("J136").Value = 860.5
$ws.Range("K136").Value = 31254003
$ws.Range("L136").Value = 2581.5
$ws.Range("M136").Value = -31251453
$ws.Range("N136").Value = -7681.5
$ws.Range("H139").Value = 70715
$ws.Range("J139").Value = 70715
$ws.Range("L139").Value = 70715
$ws.Range("N139").Value = -80995

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1750.4512
$ws.Range("I134").Value = 1075.6818
$ws.Range("J134").Value = 4533.875
$ws.Range("K134").Value = 3227.0454
$ws.Range("L134").Value = 13601.625
$ws.Range("M134").Value = -692.0454
$ws.Range("N134").Value = -18671.625
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null
$ws.Range("H139").Value = 43593.332
$ws.Range("J139").Value = 43593.332
$ws.Range("L139").Value = 43593.332
$ws.Range("N139").Value = -53873.332
$ws.Range("H140").Value = 100000
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7580702.5
$ws.Range("I31").Value = 5164.4414
$ws.Range("J31").Value = 33337532
$ws.Range("K31").Value = 5164.4414
$ws.Range("L31").Value = 33337532
$ws.Range("M31").Value = -4869.4414
$ws.Range("N31").Value = -33338122
$ws.Range("H34").Value = 7580702.5
$ws.Range("I34").Value = 5164.4414
$ws.Range("J34").Value = 33337532
$ws.Range("K34").Value = 5164.4414
$ws.Range("L34").Value = 33337532
$ws.Range("M34").Value = -4962.4414
$ws.Range("N34").Value = -33337936
$ws.Range("H132").Value = 9435230
$ws.Range("I132").Value = 10001061
$ws.Range("J132").Value = 4704.6665
$ws.Range("K132").Value = 30003183
$ws.Range("L132").Value = 14113.9995
$ws.Range("M132").Value = -30000653
$ws.Range("N132").Value = -19173.9995
$ws.Range("H134").Value = 1323.717
$ws.Range("I134").Value = 1215.449
$ws.Range("J134").Value = 2650
$ws.Range("K134").Value = 3646.347
$ws.Range("L134").Value = 7950
$ws.Range("M134").Value = -1111.347
$ws.Range("N134").Value = -13020
$ws.Range("H138").Value = 81704
$ws.Range("J138").Value = 81704
$ws.Range("L138").Value = 81704
$ws.Range("N138").Value = -91984

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 52800
$ws.Range("J138").Value = 52800
$ws.Range("L138").Value = 52800
$ws.Range("N138").Value = -63080
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
$ws.Range("H141").Value = 295743
$ws.Range("J141").Value = 295743
$ws.Range("L141").Value = 295743
$ws.Range("N141").Value = -306103

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 10317.5
$ws.Range("J104").Value = 10317.5
$ws.Range("L104").Value = 10317.5
$ws.Range("N104").Value = -17305.5
$ws.Range("H132").Value = 6762250
$ws.Range("I132").Value = 3116.131
$ws.Range("J132").Value = 38478184
$ws.Range("K132").Value = 9348.393
$ws.Range("L132").Value = 115434552
$ws.Range("M132").Value = -6818.393
$ws.Range("N132").Value = -115439612
$ws.Range("H133").Value = 64326
$ws.Range("J133").Value = 64326
$ws.Range("L133").Value = 64326
$ws.Range("N133").Value = -69386
$ws.Range("H135").Value = 139500
$ws.Range("J135").Value = 139500
$ws.Range("L135").Value = 139500
$ws.Range("N135").Value = -149640
$ws.Range("H136").Value = 11367220
$ws.Range("I136").Value = 12821689
$ws.Range("K136").Value = 38465067
$ws.Range("M136").Value = -38462517
$ws.Range("H137").Value = 80429
$ws.Range("J137").Value = 80429
$ws.Range("L137").Value = 80429
$ws.Range("N137").Value = -90629

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 53985.8
$ws.Range("J46").Value = 53985.8
$ws.Range("L46").Value = 53985.8
$ws.Range("N46").Value = -54447.8
$ws.Range("H132").Value = 1548.3334
$ws.Range("I132").Value = 1367.561
$ws.Range("J132").Value = 2118.4614
$ws.Range("K132").Value = 4102.683
$ws.Range("L132").Value = 6355.3842
$ws.Range("M132").Value = -1572.683
$ws.Range("N132").Value = -11415.3842
$ws.Range("H134").Value = 53985.8
$ws.Range("J134").Value = 53985.8
$ws.Range("L134").Value = 161957.4
$ws.Range("N134").Value = -167027.4
$ws.Range("H135").Value = 40033.08
$ws.Range("J135").Value = 40033.08
$ws.Range("L135").Value = 40033.08
$ws.Range("N135").Value = -50173.08
$ws.Range("H136").Value = 1170.7742
$ws.Range("I136").Value = 876.2143
$ws.Range("J136").Value = 3920
$ws.Range("K136").Value = 2628.6429
$ws.Range("L136").Value = 11760
$ws.Range("M136").Value = -78.64289999999983
$ws.Range("N136").Value = -16860
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null
$ws.Range("H139").Value = 41590
$ws.Range("J139").Value = 41590
$ws.Range("L139").Value = 41590
$ws.Range("N139").Value = -51870
